$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 3075.2856
$ws.Range("J97").Value = 3075.2856
$ws.Range("L97").Value = 9225.856800000001
$ws.Range("N97").Value = -10217.8568
$ws.Range("H125").Value = 1210.125
$ws.Range("I125").Value = 466
$ws.Range("J125").Value = 1458.1666
$ws.Range("K125").Value = 4194
$ws.Range("L125").Value = 13123.4994
$ws.Range("M125").Value = -1734
$ws.Range("N125").Value = -18043.4994
$ws.Range("H137").Value = 21741514
$ws.Range("I137").Value = 33335066
$ws.Range("J137").Value = 3600.625
$ws.Range("K137").Value = 100005198
$ws.Range("L137").Value = 10801.875
$ws.Range("M137").Value = -100002648
$ws.Range("N137").Value = -15901.875
$ws.Range("H138").Value = 3653
$ws.Range("I138").Value = 2647.4092
$ws.Range("J138").Value = 4267.528
$ws.Range("K138").Value = 7942.2276
$ws.Range("L138").Value = 12802.584
$ws.Range("M138").Value = -2802.2276
$ws.Range("N138").Value = -23082.584

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10638.5
$ws.Range("I32").Value = 9834.584999999999
$ws.Range("K32").Value = 9834.584999999999
$ws.Range("M32").Value = -9547.584999999999
$ws.Range("H60").Value = 200071660
$ws.Range("I60").Value = 89585.75
$ws.Range("J60").Value = 1000000000
$ws.Range("K60").Value = 89585.75
$ws.Range("L60").Value = 1000000000
$ws.Range("M60").Value = -88852.75
$ws.Range("N60").Value = -1000001466
$ws.Range("H61").Value = 14204375
$ws.Range("I61").Value = 16283228
$ws.Range("J61").Value = 1434274
$ws.Range("K61").Value = 16283228
$ws.Range("L61").Value = 1434274
$ws.Range("M61").Value = -16283016
$ws.Range("N61").Value = -1434698
$ws.Range("H97").Value = 1692.4642
$ws.Range("I97").Value = 685.8421
$ws.Range("J97").Value = 3817.5557
$ws.Range("K97").Value = 685.8421
$ws.Range("L97").Value = 3817.5557
$ws.Range("M97").Value = -189.8421
$ws.Range("N97").Value = -4809.5557
$ws.Range("H136").Value = 14204375
$ws.Range("I136").Value = 16283228
$ws.Range("J136").Value = 1434274
$ws.Range("K136").Value = 48849684
$ws.Range("L136").Value = 4302822
$ws.Range("M136").Value = -48847134
$ws.Range("N136").Value = -4307922

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 500
$ws.Range("I16").Value = 500
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 500
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -330
$ws.Range("N16").Value = ""
$ws.Range("H86").Value = 6046.6
$ws.Range("I86").Value = 4381.5713
$ws.Range("K86").Value = 4381.5713
$ws.Range("M86").Value = -3258.5713
$ws.Range("H89").Value = 6046.6
$ws.Range("I89").Value = 4381.5713
$ws.Range("K89").Value = 21907.8565
$ws.Range("M89").Value = -16291.8565

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20836098
$ws.Range("I31").Value = 30305286
$ws.Range("K31").Value = 30305286
$ws.Range("M31").Value = -30304991
$ws.Range("H34").Value = 20836098
$ws.Range("I34").Value = 30305286
$ws.Range("K34").Value = 30305286
$ws.Range("M34").Value = -30305084
$ws.Range("H58").Value = 2916.8235
$ws.Range("I58").Value = 1969.5
$ws.Range("K58").Value = 1969.5
$ws.Range("M58").Value = -1766.5
$ws.Range("H62").Value = 15166296
$ws.Range("I62").Value = 13699.556
$ws.Range("J62").Value = 25656556
$ws.Range("K62").Value = 13699.556
$ws.Range("L62").Value = 25656556
$ws.Range("M62").Value = -13075.556
$ws.Range("N62").Value = -25657804
$ws.Range("H65").Value = 15166296
$ws.Range("I65").Value = 13699.556
$ws.Range("J65").Value = 25656556
$ws.Range("K65").Value = 68497.78
$ws.Range("L65").Value = 128282780
$ws.Range("M65").Value = -65377.78
$ws.Range("N65").Value = -128289020
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").Value = ""
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").Value = ""
$ws.Range("H88").Value = 23071
$ws.Range("J88").Value = 23071
$ws.Range("L88").Value = 23071
$ws.Range("N88").Value = -23883
$ws.Range("H91").Value = 23071
$ws.Range("J91").Value = 23071
$ws.Range("L91").Value = 23071
$ws.Range("N91").Value = -25879
$ws.Range("H132").Value = 2006.3182
$ws.Range("I132").Value = 2304.75
$ws.Range("J132").Value = 1210.5
$ws.Range("K132").Value = 6914.25
$ws.Range("L132").Value = 3631.5
$ws.Range("M132").Value = -4384.25
$ws.Range("N132").Value = -8691.5
$ws.Range("H134").Value = 2008.6177
$ws.Range("I134").Value = 1776.2693
$ws.Range("K134").Value = 5328.8079
$ws.Range("M134").Value = -2793.8079
$ws.Range("H136").Value = 2916.8235
$ws.Range("I136").Value = 1969.5
$ws.Range("K136").Value = 5908.5
$ws.Range("M136").Value = -3358.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 187.36363
$ws.Range("J40").Value = 278.16666
$ws.Range("L40").Value = 1112.66664
$ws.Range("N40").Value = -1250.66664
$ws.Range("H68").Value = 397.3
$ws.Range("J68").Value = 397.3
$ws.Range("L68").Value = 1191.9
$ws.Range("N68").Value = -2813.9
$ws.Range("H71").Value = 397.3
$ws.Range("J71").Value = 397.3
$ws.Range("L71").Value = 3575.7
$ws.Range("N71").Value = -11687.7
$ws.Range("H80").Value = 55558224
$ws.Range("I80").Value = 83334830
$ws.Range("K80").Value = 250004490
$ws.Range("M80").Value = -250003554
$ws.Range("H83").Value = 55558224
$ws.Range("I83").Value = 83334830
$ws.Range("K83").Value = 750013470
$ws.Range("M83").Value = -750008790
$ws.Range("H119").Value = 8569.888999999999
$ws.Range("I119").Value = 3984.5
$ws.Range("J119").Value = 9880
$ws.Range("K119").Value = 11953.5
$ws.Range("L119").Value = 29640
$ws.Range("M119").Value = -7115.5
$ws.Range("N119").Value = -39316
$ws.Range("H121").Value = 2459.2104
$ws.Range("I121").Value = 554.4
$ws.Range("K121").Value = 1663.2
$ws.Range("M121").Value = -353.1999999999998
$ws.Range("H129").Value = 10873911
$ws.Range("I129").Value = 22730022
$ws.Range("J129").Value = 5810.5835
$ws.Range("K129").Value = 68190066
$ws.Range("L129").Value = 17431.7505
$ws.Range("M129").Value = -68185066
$ws.Range("N129").Value = -27431.7505
$ws.Range("H138").Value = 10924.7
$ws.Range("I138").Value = 13196.272
$ws.Range("K138").Value = 39588.81600000001
$ws.Range("M138").Value = -34448.81600000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 9083.666999999999
$ws.Range("I22").Value = 10898.4
$ws.Range("J22").Value = 10
$ws.Range("K22").Value = 10898.4
$ws.Range("L22").Value = 10
$ws.Range("M22").Value = -10369.4
$ws.Range("N22").Value = -1068
$ws.Range("H58").Value = 34955.184
$ws.Range("I58").Value = 26612
$ws.Range("J58").Value = 72499.5
$ws.Range("K58").Value = 26612
$ws.Range("L58").Value = 72499.5
$ws.Range("M58").Value = -26335
$ws.Range("N58").Value = -73053.5
$ws.Range("H70").Value = 14345.242
$ws.Range("J70").Value = 15086.523
$ws.Range("L70").Value = 15086.523
$ws.Range("N70").Value = -15626.523
$ws.Range("H73").Value = 14345.242
$ws.Range("J73").Value = 15086.523
$ws.Range("L73").Value = 15086.523
$ws.Range("N73").Value = -16958.523
$ws.Range("H122").Value = 7840.476
$ws.Range("I122").Value = 5884.4375
$ws.Range("K122").Value = 17653.3125
$ws.Range("M122").Value = -15203.3125

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H88").Value = 59999.5
$ws.Range("J88").Value = 59999.5
$ws.Range("L88").Value = 59999.5
$ws.Range("N88").Value = -60855.5
$ws.Range("H91").Value = 59999.5
$ws.Range("J91").Value = 59999.5
$ws.Range("L91").Value = 59999.5
$ws.Range("N91").Value = -62963.5
$ws.Range("H100").Value = 19252724
$ws.Range("I100").Value = 2546.3
$ws.Range("K100").Value = 2546.3
$ws.Range("M100").Value = -2005.3
$ws.Range("H116").Value = 119999.5
$ws.Range("J116").Value = 119999.5
$ws.Range("L116").Value = 119999.5
$ws.Range("N116").Value = -129177.5
$ws.Range("H128").Value = 66283.664
$ws.Range("J128").Value = 66283.664
$ws.Range("L128").Value = 66283.664
$ws.Range("N128").Value = -76243.664

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1470.2858
$ws.Range("I81").Value = 1382
$ws.Range("K81").Value = 2764
$ws.Range("M81").Value = -1703
$ws.Range("H84").Value = 1470.2858
$ws.Range("I84").Value = 1382
$ws.Range("K84").Value = 13820
$ws.Range("M84").Value = -8516
$ws.Range("H100").Value = 1036.5483
$ws.Range("I100").Value = 1010.13635
$ws.Range("K100").Value = 2020.2727
$ws.Range("M100").Value = -1479.2727
$ws.Range("H105").Value = 34998
$ws.Range("J105").Value = 34998
$ws.Range("L105").Value = 34998
$ws.Range("N105").Value = -41986
$ws.Range("H122").Value = 2220.2273
$ws.Range("I122").Value = 1622.4117
$ws.Range("K122").Value = 4867.2351
$ws.Range("M122").Value = -2417.2351
$ws.Range("H140").Value = 44475.2
$ws.Range("J140").Value = 44475.2
$ws.Range("L140").Value = 44475.2
$ws.Range("N140").Value = -54835.2
$ws.Range("H141").Value = 89182.89
$ws.Range("J141").Value = 89182.89
$ws.Range("L141").Value = 89182.89
$ws.Range("N141").Value = -99542.89
